$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.658.70"
$ws.Range("D3").Value = "'1.597.56"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'211.46"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'19.57"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'1.821.89"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'1.592.14"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'65.06"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'26.651.61"
$ws.Range("E18").Value = "  +1.37%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").Value = "'209.08"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'7.02"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'2.32"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'8.99"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'144.15"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'7.13"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'0.0515"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "'1.288.78"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").Value = "'0.616"
$ws.Range("E35").Value = "  -7.53%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +17.86%  "
$ws.Range("D41").Value = "'5.50"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "'63.54"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'1.734.09"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'90.71"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "'7.42"
$ws.Range("E51").Value = "  -1.08%  "
